$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Template text updates (JXLS expressions) ---

# Period: now uses from/to .toString(...) with explicit date pattern instead of
# the old "".format(...) construct.
$ws.Range("B6").Value = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'

# position.fixTime now converted through joda DateTime with the client timezone
$ws.Range("B9").Value = '${new("org.joda.time.DateTime", position.fixTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'

# OpenStreetMap link now uses https instead of http
$ws.Range("G9").Value = '${util.hyperlink("".format("https://www.openstreetmap.org/?mlat=%1$f&mlon=%2$f#map=16/%1$f/%2$f", position.latitude, position.longitude), position.getAddress() == null ? "".format("%1$f°, %2$f°", position.latitude, position.longitude) : position.address)}'

# --- Formatting tweaks ---

# Indent levels bumped up slightly on the label/value columns
$ws.Range("B1").IndentLevel = 15
$ws.Range("B3").IndentLevel = 15

$ws.Range("B2").IndentLevel = 2
$ws.Range("B4").IndentLevel = 2
$ws.Range("B5").IndentLevel = 2
$ws.Range("B6").IndentLevel = 2

# Slightly wider G/H columns
$ws.Columns.Item(7).ColumnWidth = 61.85546875
$ws.Columns.Item(8).ColumnWidth = 73.28515625

# Update the selected cell shown when the sheet is opened
$ws.Range("G9").Select()
